# Generate Report for Handback
# Update the handback timestamps / priority status for the
# 46c8fcce-46a7-42bb-bb59-22347e7128cc.md file across the Overview,
# zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 46c8fcce row
# (row 3) and the afadae09 row (row 4) which shared the same old value.
$wsOverview.Range("G3").Value = "2016-08-27 16:17:43"
$wsOverview.Range("G4").Value = "2016-08-27 16:17:43"

# zh-cn sheet: Priority ("ht" -> "mt") for rows 3 and 4
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime for rows 3 and 4
$wsZhCn.Range("H3").Value = "2016-08-27 16:17:39"
$wsZhCn.Range("H4").Value = "2016-08-27 16:17:39"

# zh-cn sheet: Correspond Handback DateTime for rows 3 and 4
$wsZhCn.Range("K3").Value = "2016-08-27 16:17:58"
$wsZhCn.Range("K4").Value = "2016-08-27 16:17:58"

# de-de sheet: Priority ("ht" -> "mt") for rows 3 and 4
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# de-de sheet: Correspond Handoff Datetime for rows 3 and 4
$wsDeDe.Range("H3").Value = "2016-08-27 16:17:43"
$wsDeDe.Range("H4").Value = "2016-08-27 16:17:43"

# de-de sheet: Correspond Handback DateTime for rows 3 and 4
$wsDeDe.Range("K3").Value = "2016-08-27 16:18:10"
$wsDeDe.Range("K4").Value = "2016-08-27 16:18:10"
